# Aggiornamento Variable e AnalysisUnit per BE
# Adds 16 new rows (5-20) to the "r AnalysisUnit_Variable" sheet with the
# new BE indicator variables.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("r AnalysisUnit_Variable")

# Data for the new rows: row number -> (B/C value, F value)
# Row 7 (IND_3) is entered first, then the remaining B values top-to-bottom,
# then the remaining F values top-to-bottom, matching how the rows were
# actually authored.
$row7Code = "CUSTOMER_BE_IND_3"
$row7Target = "BE_IND_3"

$otherRows = @(
    @{ Row = 5;  Code = "CUSTOMER_BE_IND_1";  Target = "BE_IND_1"  },
    @{ Row = 6;  Code = "CUSTOMER_BE_IND_2";  Target = "BE_IND_2"  },
    @{ Row = 8;  Code = "CUSTOMER_BE_IND_7";  Target = "BE_IND_7"  },
    @{ Row = 9;  Code = "CUSTOMER_BE_IND_8";  Target = "BE_IND_8"  },
    @{ Row = 10; Code = "CUSTOMER_BE_IND_9";  Target = "BE_IND_9"  },
    @{ Row = 11; Code = "CUSTOMER_BE_IND_14"; Target = "BE_IND_14" },
    @{ Row = 12; Code = "CUSTOMER_BE_IND_16"; Target = "BE_IND_16" },
    @{ Row = 13; Code = "CUSTOMER_BE_IND_34"; Target = "BE_IND_34" },
    @{ Row = 14; Code = "CUSTOMER_BE_IND_35"; Target = "BE_IND_35" },
    @{ Row = 15; Code = "CUSTOMER_BE_IND_40"; Target = "BE_IND_40" },
    @{ Row = 16; Code = "CUSTOMER_BE_IND_44"; Target = "BE_IND_44" },
    @{ Row = 17; Code = "CUSTOMER_BE_IND_48"; Target = "BE_IND_48" },
    @{ Row = 18; Code = "CUSTOMER_BE_IND_51"; Target = "BE_IND_51" },
    @{ Row = 19; Code = "CUSTOMER_BE_IND_55"; Target = "BE_IND_55" },
    @{ Row = 20; Code = "CUSTOMER_BE_IND_56"; Target = "BE_IND_56" }
)

# 1) Row 7 first (the "IND_3" row), A, B, C, E, F all together.
$ws.Cells.Item(7, 1).Value = "CREATE/MODIFY"
$ws.Cells.Item(7, 2).Value = $row7Code
$ws.Cells.Item(7, 3).Value = $row7Code
$ws.Cells.Item(7, 5).Value = "COUNTERPARTY_BIB"
$ws.Cells.Item(7, 6).Value = $row7Target

# 2) Column A for the remaining rows.
foreach ($item in $otherRows) {
    $ws.Cells.Item($item.Row, 1).Value = "CREATE/MODIFY"
}

# 3) Column B for the remaining rows (top to bottom).
foreach ($item in $otherRows) {
    $ws.Cells.Item($item.Row, 2).Value = $item.Code
}

# 4) Column C for the remaining rows (top to bottom) - mirrors column B.
foreach ($item in $otherRows) {
    $ws.Cells.Item($item.Row, 3).Value = $item.Code
}

# 5) Column E for the remaining rows (top to bottom).
foreach ($item in $otherRows) {
    $ws.Cells.Item($item.Row, 5).Value = "COUNTERPARTY_BIB"
}

# 6) Column F for the remaining rows (top to bottom).
foreach ($item in $otherRows) {
    $ws.Cells.Item($item.Row, 6).Value = $item.Target
}

# Update the selection state to reflect the last edited cells (rows 11,
# 14, 20 and 21 highlighted, active cell at A21).
$ws.Rows.Item(11).Select() | Out-Null
$ws.Rows.Item(14).Select() | Out-Null
$ws.Rows.Item(20).Select() | Out-Null
$ws.Rows.Item(21).Select() | Out-Null

Write-Output "Sheet 'r AnalysisUnit_Variable' updated with BE indicator rows 5-20."
